$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a brand-new row at row 17. This pushes the old rows 17-26
#    ("Min:" .. "Forth Q:") down to rows 18-27 while leaving rows 14-16
#    ("Median:", "Mode:", "StDev:") exactly where they are, matching the
#    target layout (Var.S is the newly inserted stat row).
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Insert()

# ---------------------------------------------------------------------------
# 2. New "string function" columns (A/B/C) next to the existing stats table.
# ---------------------------------------------------------------------------

# Row 14: Pi label + PI() formula, plus UPPER() of the Yellow cell (C8).
$ws.Range("A14").Value = "Pi:"
$ws.Range("A14").HorizontalAlignment = -4152
$ws.Range("B14").Formula = "=PI()"
$ws.Range("B14").NumberFormat = "0.000"
$ws.Range("C14").Formula = "=UPPER(C8)"

# Row 15: blank placeholder under Pi, plus LOWER() of the named range Cyan.
$ws.Range("B15").NumberFormat = "0.000"
$ws.Range("C15").Formula = "=LOWER(Cyan)"

# Row 16: Yellow Len label + LEN() formula.
$ws.Range("B16").Value = "Yellow Len:"
$ws.Range("B16").NumberFormat = "0.000"
$ws.Range("C16").Formula = "=LEN(C8)"

# Row 17 (new row): Trimmed label + TRIM(CONCATENATE()) formula, plus the
# new Var (S) stat in the D/E columns.
$ws.Range("B17").Value = "Trimmed:"
$ws.Range("B17").NumberFormat = "0.000"
$ws.Range("C17").Formula = '=TRIM(CONCATENATE(C8, "   "))'
$ws.Range("D17").Value = "Var (S)"
$ws.Range("D17").HorizontalAlignment = -4152
$ws.Range("E17").Formula = "=VAR(E2:E9)"

# ---------------------------------------------------------------------------
# 3. Rename the StDev label (row 16) to the new "StDev (S):" text.
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = "StDev (S):"

# ---------------------------------------------------------------------------
# 4. Defined names: drop the two that were removed, shift/rename the ones
#    whose target rows moved, and add the brand-new ones.
# ---------------------------------------------------------------------------
$wb.Names.Item("Booleans").Delete()
$wb.Names.Item("Subset").Delete()

$wb.Names.Item("DevSq").RefersTo    = '=''Sample 2''!$E$23'
$wb.Names.Item("FirstQ").RefersTo   = '=''Sample 2''!$E$24'
$wb.Names.Item("ForthQ").RefersTo   = '=''Sample 2''!$E$27'
$wb.Names.Item("Kurtosis").RefersTo = '=''Sample 2''!$E$22'
$wb.Names.Item("Max").RefersTo      = '=''Sample 2''!$E$19'
$wb.Names.Item("Min").RefersTo      = '=''Sample 2''!$E$18'
$wb.Names.Item("SecondQ").RefersTo  = '=''Sample 2''!$E$25'
$wb.Names.Item("Skew").RefersTo     = '=''Sample 2''!$E$20'
$wb.Names.Item("SumSq").RefersTo    = '=''Sample 2''!$E$21'
$wb.Names.Item("ThirdQ").RefersTo   = '=''Sample 2''!$E$26'

$wb.Names.Item("StDev").Name = "StDev.S"

$wb.Names.Add("LOWERCASE", '=''Sample 2''!$C$15')
$wb.Names.Add("Pi", '=''Sample 2''!$B$14')
$wb.Names.Add("UPPERCASE", '=''Sample 2''!$C$14')
$wb.Names.Add("Var.S", '=''Sample 2''!$E$17')
$wb.Names.Add("YellowLen", '=''Sample 2''!$C$16')
$wb.Names.Add("YellowTrim", '=''Sample 2''!$C$17')

# ---------------------------------------------------------------------------
# 5. Workbook window geometry update (best-effort; matches the saved
#    workbookView window position/size).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 7800
$win.Top = 320
$win.Width = 14440
$win.Height = 16060
